$d = $word.ActiveDocument

# Use tracked changes while editing so the engine keeps newly-typed text
# as its own run instead of silently re-coalescing it into a neighbouring
# run that happens to share the same formatting (that's what real Word
# does too: a fresh edit gets its own <w:r>, merge only happens when you
# explicitly ask for it).
$d.TrackRevisions = $true

# ---------------------------------------------------------------------
# 1) "...${periodoConvenio} entre el Servicio de Salud Iquique..."
#    -> "...${periodoConvenio}, entre el Servicio de Salud Iquique..."
#    The leading space of the " entre el Servicio de Salud " run becomes
#    ", " (comma + space), splitting that run in two.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" entre el Servicio de Salud", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$lead = $d.Range($rng.Start, $rng.Start + 1)
$lead.Text = ", "

# ---------------------------------------------------------------------
# 2) "...al ítem 24-03-298-002..." -> "...al ítem N°24-03 298-002..."
#    a) trailing space of "al ítem " becomes " N°" (splitting the run)
#    b) the "-" run between "24-03" and "298-002" becomes " "
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("al ítem ", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$trail = $d.Range($rng.End - 1, $rng.End)
$trail.Text = " N°"

$rng = $d.Content
$rng.Find.Execute("24-03-298-002", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$hyphen = $d.Range($rng.Start + 5, $rng.Start + 6)
$hyphen.Text = " "

$d.TrackRevisions = $false
$d.AcceptAllRevisions()
